# Add two new summary-annotation text boxes to slide 1, as introduced by the
# "new simulations and summary files" commit:
#   - TextBox 2  (id 3)  - three lines of "~(p=.., cl=.., np=..)" text
#   - TextBox 11 (id 12) - two lines naming the new Sim119 / Sim120 runs
#
# PowerPoint's COM Shapes.AddTextbox/Left/Top/Width/Height are expressed in
# points, while the target OOXML <a:off>/<a:ext> values are EMU, so convert
# (1 pt = 12700 EMU) to land on the exact EMU coordinates from the diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function New-AnnotationTextBox($offX, $offY, $extCx, $extCy, [string[]]$lines) {
    $left   = $offX / 12700.0
    $top    = $offY / 12700.0
    $width  = $extCx / 12700.0
    $height = $extCy / 12700.0

    $tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
    $tb.TextFrame.WordWrap = -1
    $tb.TextFrame.AutoSize = 1
    $tb.Fill.Visible = 0

    for ($i = 0; $i -lt $lines.Count; $i++) {
        if ($i -eq 0) {
            $tr = $tb.TextFrame.TextRange
            $tr.Text = $lines[$i]
        } else {
            [void]$tb.TextFrame.TextRange.InsertAfter("`r" + $lines[$i])
            $tr = $tb.TextFrame.TextRange.Lines($i + 1)
        }
        $tr.Font.Size = 22
        $tr.LanguageID = "en-AU"
    }

    return $tb
}

[void](New-AnnotationTextBox 9622653 3739654 2495365 1107996 @(
    "~(p=8, cl=3, np=5)",
    "~(p=8, cl=4, np=4)",
    "~(p=10, cl=4, np=5)"
))

[void](New-AnnotationTextBox 8842160 5589544 3124544 769441 @(
    "Sim119~(p=6, cl=3, np=2)",
    "Sim120~(p=6, cl=3, np=2)"
))
